$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new paragraph right after the Heading1 title paragraph
#    containing a bold "Meta description" run followed by the
#    description text run.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range

$metaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Find out why Berryburst Max is one of the hottest new slot games. Enjoy engaging graphics and intense gameplay with bonuses and promotions available.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaInsertResult = $metaRange.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Berryburst Max for Free - Review
#    and Bonuses" paragraph that used to sit right before the final
#    italic meta-description paragraph at the end of the document.
#    (Paragraph 1 carries the same text as the Heading1 title and must
#    be left untouched - only the later, duplicated copy goes away.)
# ---------------------------------------------------------------------
$target = $null
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.Trim()
    if (($txt -eq "Play Berryburst Max for Free - Review and Bonuses") -and ($i -ne 1)) {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# ---------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    feature-image prompt, keeping its leading empty run + italic
#    formatting intact.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$textEnd = $lastRange.End - 1
$replaceRange = $d.Range($lastRange.Start, $textEnd)

$newImageXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Please create a cartoon-style feature image that portrays a happy Maya warrior with glasses for the game &quot;Berryburst Max&quot;. The image should feature the Maya warrior as the main focus, with bold colors and bright, fruity elements surrounding them. The Maya warrior should be holding a slot machine lever with one hand, and a big smile on their face indicating a big win. The glasses should be thick-framed and add to the playful and colorful nature of the image. Please feel free to add any other elements that fit the theme and tone of the game, such as stylized fruit symbols or sparkles representing the excitement of winning. The overall goal of the feature image is to capture the fun and energetic nature of the game and encourage players to take a spin.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$imageInsertResult = $replaceRange.InsertXML($newImageXml)

Write-Host "Edit complete"
